$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.470.39"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.852.75"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6315"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07680"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07758"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "1.861.08"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.041"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6814"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001070"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "2.110.19"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.174"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "29.488.26"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.453"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1385"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.411"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.324"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.473"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05690"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.138"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.850"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.167"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7096"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.590"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.784"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01796"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").Value = "1.223.18"
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.558"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.74%  "
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.142"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4026"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.010"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.686"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1142"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.31%  "
